$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.398627400398254
$ws.Range("B1").Value = 2.596854209899902
$ws.Range("C1").Value = 6.710230350494385
$ws.Range("D1").Value = 2.418612718582153
$ws.Range("E1").Value = 1.195956945419312
